# Add a DB25 connector (X2) line to the partlist BOM.
#
# The new part is inserted as a whole new row right before the existing
# "shunt jumpers" row (previously row 34, now shifted to row 35), so the
# row order becomes:
#   ...
#   33  ET-3400 header
#   34  X2 / DB25 / 1 / DB25 connector / Digikey, AE10935-ND   <-- new
#   35  (blank Part) / shunt jumpers / ... / Digikey, S9337-ND
#   36  CAB1 / 40-ribbon / ...
#   37  IDC1, IDC2 / 2x20 IDC conn / ...

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 34, pushing everything from the old row 34
# down by one (old 34->35, 35->36, 36->37).
$ws.Rows.Item(34).Insert()

# Populate the new row: Part, Value, Quantity, Description, Source.
$ws.Cells.Item(34, 1).Value = "X2"
$ws.Cells.Item(34, 2).Value = "DB25"
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = "DB25 connector"
$ws.Cells.Item(34, 5).Value = "Digikey, AE10935-ND"

# Match the author's final selection/active cell in the sheet.
$ws.Range("G43").Select() | Out-Null
